$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.293.71"
$ws.Range("E2").Value = "  +0.98%  "

# Row 3
$ws.Range("D3").Value = "2.247.09"
$ws.Range("E3").Value = "  +0.46%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.67%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.573"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.24%  "

# Row 8
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.89%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.79%  "

# Row 11
$ws.Range("E11").Value = "  -0.96%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.23%  "

# Row 13
$ws.Range("E13").Value = "  +0.05%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.590.82"
$ws.Range("E14").Value = "  +0.51%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.332.52"
$ws.Range("E15").Value = "  +4.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.833"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.94%  "

# Row 18
$ws.Range("D18").Value = "44.079.46"
$ws.Range("E18").Value = "  +0.83%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0965"
$ws.Range("E19").Value = "  -0.39%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.72%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.64%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.42%  "

# Row 26
$ws.Range("E26").Value = "  +0.18%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.28%  "

# Row 28
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.25%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.51%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.62%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.11%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.85%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0797"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.77%  "

# Row 34
$ws.Range("E34").Value = "  -0.63%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.04%  "

# Row 36
$ws.Range("E36").Value = "  +2.36%  "

# Row 37
$ws.Range("E37").Value = "  -0.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.37%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.72%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.35%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0301"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.74%  "

# Row 43
$ws.Range("E43").Value = "  +0.24%  "

# Row 44
$ws.Range("D44").Value = "1.748.74"
$ws.Range("E44").Value = "  +2.88%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "82.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.192"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.18%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.86%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.14%  "

# Row 50
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.05%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.28%  "
